$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Collect all D-column cell addresses that receive new values so we can
# temporarily mark them as Text before assignment (stops Excel's automatic
# number/date coercion of values like '581.68' or '19.03'), then restore
# the original (default) cell style afterwards so no stray style index is
# left behind on the cell.
$dCells = @("D2","D3","D5","D6","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D25","D26","D27","D28","D29","D30","D31","D32","D35","D36","D41","D42","D43","D46","D47","D48","D49","D50","D51","D38","D39","D44","D45")
foreach ($addr in $dCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = "73.011.73"
$ws.Range("E2").Value = "  +7.59%  "
$ws.Range("D3").Value = "2.567.53"
$ws.Range("E3").Value = "  +6.85%  "
$ws.Range("E4").Value = "  -0.13%  "
$ws.Range("D5").Value = "581.68"
$ws.Range("E5").Value = "  +5.11%  "
$ws.Range("D6").Value = "179.43"
$ws.Range("E6").Value = "  +13.38%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "0.529"
$ws.Range("E8").Value = "  +4.83%  "
$ws.Range("D9").Value = "0.197"
$ws.Range("E9").Value = "  +22.26%  "
$ws.Range("D10").Value = "2.564.32"
$ws.Range("E10").Value = "  +6.78%  "
$ws.Range("D11").Value = "0.162"
$ws.Range("E11").Value = "  -0.66%  "
$ws.Range("D12").Value = "0.351"
$ws.Range("E12").Value = "  +6.69%  "
$ws.Range("D13").Value = "4.79"
$ws.Range("E13").Value = "  +3.25%  "
$ws.Range("D14").Value = "0.0000191"
$ws.Range("E14").Value = "  +10.76%  "
$ws.Range("D15").Value = "72.932.13"
$ws.Range("E15").Value = "  +7.62%  "
$ws.Range("D16").Value = "2.992.82"
$ws.Range("E16").Value = "  +5.18%  "
$ws.Range("D17").Value = "25.55"
$ws.Range("E17").Value = "  +12.34%  "
$ws.Range("D18").Value = "2.556.74"
$ws.Range("E18").Value = "  +6.19%  "
$ws.Range("D19").Value = "11.49"
$ws.Range("E19").Value = "  +11.48%  "
$ws.Range("D20").Value = "7.79"
$ws.Range("E20").Value = "  +14.16%  "
$ws.Range("D21").Value = "360.46"
$ws.Range("E21").Value = "  +9.44%  "
$ws.Range("D22").Value = "2.21"
$ws.Range("E22").Value = "  +19.33%  "
$ws.Range("D23").Value = "4.02"
$ws.Range("E23").Value = "  +6.52%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "68.71"
$ws.Range("E25").Value = "  +4.25%  "
$ws.Range("D26").Value = "4.04"
$ws.Range("E26").Value = "  +11.45%  "
$ws.Range("D27").Value = "8.99"
$ws.Range("E27").Value = "  +10.78%  "
$ws.Range("D28").Value = "2.674.53"
$ws.Range("E28").Value = "  +5.49%  "
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.33%  "
$ws.Range("D30").Value = "0.0₃0934"
$ws.Range("E30").Value = "  +16.22%  "
$ws.Range("D31").Value = "506.24"
$ws.Range("E31").Value = "  +21.52%  "
$ws.Range("D32").Value = "7.77"
$ws.Range("E32").Value = "  +10.18%  "
$ws.Range("E33").Value = "  +16.01%  "
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("D35").Value = "1.69"
$ws.Range("E35").Value = "  +6.44%  "
$ws.Range("D36").Value = "158.71"
$ws.Range("E36").Value = "  -0.34%  "
$ws.Range("E37").Value = "  +10.50%  "
$ws.Range("D41").Value = "4.79"
$ws.Range("E41").Value = "  +12.20%  "
$ws.Range("D42").Value = "1.63"
$ws.Range("E42").Value = "  +11.19%  "
$ws.Range("D43").Value = "0.316"
$ws.Range("E43").Value = "  +7.42%  "
$ws.Range("D46").Value = "1.14"
$ws.Range("E46").Value = "  +7.40%  "
$ws.Range("D47").Value = "147.84"
$ws.Range("E47").Value = "  +12.87%  "
$ws.Range("D48").Value = "3.54"
$ws.Range("E48").Value = "  +7.68%  "
$ws.Range("D49").Value = "0.515"
$ws.Range("E49").Value = "  +8.29%  "
$ws.Range("D50").Value = "0.0753"
$ws.Range("E50").Value = "  +6.37%  "
$ws.Range("D51").Value = "0.583"
$ws.Range("E51").Value = "  +5.81%  "

$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").Value = "19.03"
$ws.Range("E38").Value = "  +7.39%  "
$ws.Range("B39").Value = "WhiteBITCoin"
$ws.Range("C39").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D39").Value = "19.33"
$ws.Range("E39").Value = "  +1.85%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.36"
$ws.Range("E44").Value = "  +21.09%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "38.48"
$ws.Range("E45").Value = "  +3.11%  "

# Restore default (unstyled) formatting on the D cells we touched above.
foreach ($addr in $dCells) { $ws.Range($addr).Style = "Normal" }
